# Documentation Checklist update
# Mark "Alvin" as the documenter (column B) for a few more source files
# that have since been documented/checked.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B15").Value2 = "Alvin"
$ws.Range("B38").Value2 = "Alvin"
$ws.Range("B69").Value2 = "Alvin"
$ws.Range("B70").Value2 = "Alvin"

# Scroll the view down to where the latest edits were made.
$ws.Application.Goto($ws.Range("A51"), $false)
$ws.Range("B70").Select() | Out-Null
